# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that hold the "Component 3" tables) get
#    switched from the custom "Table_0" style to the built-in
#    "No Style, Table Grid" style.
# 2) The deck's theme colour scheme (the part actually driving the slide
#    master / slides, ppt/theme/theme2.xml) is switched from the
#    "Red Violet" ("Integral") palette back to the stock "Office" palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Re-style the three tables.
# ---------------------------------------------------------------------
$newTableStyleId = "{28A34EEA-2248-4B0C-876F-9FC12AC68CE6}"

for ($slideIdx = 14; $slideIdx -le 16; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2. Swap the theme colour palette back to the "Office" scheme.
# ---------------------------------------------------------------------
function ConvertTo-OleColor($rrggbb) {
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme.Item(1..12):
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $colorScheme.Item($i + 1).RGB = ConvertTo-OleColor($officeColors[$i])
}
